# Expand example biobank gender reference data (LifeLines_Gender and
# Prevend_Gender sheets) with extra gender categories, used for testing
# workbooks with more than 10 categories.

$wb = $excel.ActiveWorkbook

$newGenders = @("Vulcan", "Krogan", "Asari", "Martian", "Jupitaan", "American ", "German", "Venusaur", "Charizard ")

# --- LifeLines_Gender --------------------------------------------------
$llg = $wb.Worksheets.Item("LifeLines_Gender")

# existing Male/Female rows pick up the plain-black-font style already
# used elsewhere in the workbook (style index 4)
$llg.Range("A2:B3").Font.Color = 0

$row = 4
$value = 2
foreach ($gender in $newGenders) {
    $llg.Cells.Item($row, 1).Value = $gender
    $llg.Cells.Item($row, 2).Value = $value
    $row = $row + 1
    $value = $value + 1
}
$llg.PageSetup.PaperSize = 9
$llg.PageSetup.Orientation = 1

# --- Prevend_Gender ------------------------------------------------------
$pg = $wb.Worksheets.Item("Prevend_Gender")

$row = 4
$value = 3
foreach ($gender in $newGenders) {
    $pg.Cells.Item($row, 1).Value = $gender
    $pg.Cells.Item($row, 2).Value = $value
    $row = $row + 1
    $value = $value + 1
}
$pg.PageSetup.PaperSize = 9
$pg.PageSetup.Orientation = 1
$pg.Range("A2:B12").Select()

# --- view state: LifeLines_Gender becomes the active tab/selection -----
$attrs = $wb.Worksheets.Item("attributes")
$attrs.Range("D26").Select()

$llg.Activate()
$llg.Range("E6").Select()
